$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 603 ("「うちゅうへとびたい」...") entirely; all rows below shift up by one.
$ws.Rows.Item(603).Delete()
